$wb = $excel.ActiveWorkbook

# --- Installation sheet: drop the "Electric_boiler" column (old column B) ---
$wsInstall = $wb.Worksheets.Item("Installation")
$wsInstall.Range("B1").EntireColumn.Delete()

# --- Capacity sheet: drop the "Electric_boiler" column (old column B) ---
$wsCapacity = $wb.Worksheets.Item("Capacity")
$wsCapacity.Range("B1").EntireColumn.Delete()

# Capacity sheet values were also updated (price policy re-run), not just shifted.
$wsCapacity.Range("B2").Value = 69.096061656646114
$wsCapacity.Range("F2").Value = 2339.0594869023944
$wsCapacity.Range("B3").Value = 119.53618666599779
$wsCapacity.Range("C3").Value = 459.16087512338322

# --- Storage_capacity sheet: updated total ---
$wsStorage = $wb.Worksheets.Item("Storage_capacity")
$wsStorage.Range("B2").Value = 959.35649950646734
